# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K") holds per-start strikeout counts. This regenerates
# those values (rows 2-36) with the newly-calculated "K" figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G, rows 2 through 36 (in row order).
$kVals = @(10, 3, 3, 1, 8, 7, 2, 2, 6, 4, 4, 5, 5, 3, 5, 3, 6, 4, 7, 8, 6, 2, 6, 3, 2, 9, 11, 5, 3, 4, 7, 4, 1, 0, 3)

$startRow = 2
for ($i = 0; $i -lt $kVals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kVals[$i]
}
